$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 903.4
$ws.Range("I28").Value = 837.1111
$ws.Range("K28").Value = 837.1111
$ws.Range("M28").Value = -352.1111
$ws.Range("H86").Value = 2590.9412
$ws.Range("J86").Value = 3460.1
$ws.Range("L86").Value = 3460.1
$ws.Range("N86").Value = -5706.1
$ws.Range("H88").Value = 2035.3889
$ws.Range("I88").Value = 1778.6666
$ws.Range("J88").Value = 2086.7334
$ws.Range("K88").Value = 1778.6666
$ws.Range("L88").Value = 2086.7334
$ws.Range("M88").Value = -1372.6666
$ws.Range("N88").Value = -2898.7334
$ws.Range("H89").Value = 2590.9412
$ws.Range("J89").Value = 3460.1
$ws.Range("L89").Value = 17300.5
$ws.Range("N89").Value = -28532.5
$ws.Range("H91").Value = 2035.3889
$ws.Range("I91").Value = 1778.6666
$ws.Range("J91").Value = 2086.7334
$ws.Range("K91").Value = 1778.6666
$ws.Range("L91").Value = 2086.7334
$ws.Range("M91").Value = -374.6666
$ws.Range("N91").Value = -4894.7334
$ws.Range("H97").Value = 2244.8333
$ws.Range("J97").Value = 2824
$ws.Range("L97").Value = 8472
$ws.Range("N97").Value = -9464
$ws.Range("H99").Value = 1279.4166
$ws.Range("I99").Value = 1647.625
$ws.Range("J99").Value = 543
$ws.Range("K99").Value = 4942.875
$ws.Range("L99").Value = 1629
$ws.Range("M99").Value = -3444.875
$ws.Range("N99").Value = -4625
$ws.Range("H100").Value = 6982.2856
$ws.Range("I100").Value = 1556.8334
$ws.Range("J100").Value = 11051.375
$ws.Range("K100").Value = 1556.8334
$ws.Range("L100").Value = 11051.375
$ws.Range("M100").Value = -1015.8334
$ws.Range("N100").Value = -12133.375
$ws.Range("H101").Value = 494.91666
$ws.Range("I101").Value = 509
$ws.Range("K101").Value = 1527
$ws.Range("M101").Value = 95
$ws.Range("H106").Value = 4355.875
$ws.Range("I106").Value = 4999
$ws.Range("J106").Value = 4264
$ws.Range("K106").Value = 4999
$ws.Range("L106").Value = 4264
$ws.Range("M106").Value = -4368
$ws.Range("N106").Value = -5526
$ws.Range("H107").Value = 943.13635
$ws.Range("I107").Value = 1047.4286
$ws.Range("J107").Value = 760.625
$ws.Range("K107").Value = 1047.4286
$ws.Range("L107").Value = 760.625
$ws.Range("M107").Value = 872.5714
$ws.Range("N107").Value = -4600.625
$ws.Range("H116").Value = 3049.4546
$ws.Range("I116").Value = 2788.889
$ws.Range("K116").Value = 2788.889
$ws.Range("M116").Value = 653.1109999999999
$ws.Range("H138").Value = 177732.98
$ws.Range("J138").Value = 252878.38
$ws.Range("L138").Value = 758635.14
$ws.Range("N138").Value = -768915.14

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4408.6665
$ws.Range("I132").Value = 4406.0713
$ws.Range("K132").Value = 13218.2139
$ws.Range("M132").Value = -10688.2139

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3778.5454
$ws.Range("I105").Value = 3530.0303
$ws.Range("J105").Value = 4524.091
$ws.Range("K105").Value = 3530.0303
$ws.Range("L105").Value = 4524.091
$ws.Range("M105").Value = -1783.0303
$ws.Range("N105").Value = -8018.091

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2567.5
$ws.Range("I58").Value = 2625.4285
$ws.Range("K58").Value = 2625.4285
$ws.Range("M58").Value = -2422.4285
$ws.Range("H105").Value = 1214.6666
$ws.Range("I105").Value = 1322.75
$ws.Range("K105").Value = 1322.75
$ws.Range("M105").Value = 424.25
$ws.Range("H116").Value = 150000
$ws.Range("J116").Value = 150000
$ws.Range("L116").Value = 150000
$ws.Range("N116").Value = -159178
$ws.Range("H134").Value = 4239.476
$ws.Range("I134").Value = 1877.4166
$ws.Range("K134").Value = 5632.2498
$ws.Range("M134").Value = -3097.2498
$ws.Range("H136").Value = 2567.5
$ws.Range("I136").Value = 2625.4285
$ws.Range("K136").Value = 7876.2855
$ws.Range("M136").Value = -5326.2855

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8830535
$ws.Range("I4").Value = 7876605
$ws.Range("K4").Value = 23629815
$ws.Range("M4").Value = -23629703
$ws.Range("H96").Value = 264981.75
$ws.Range("J96").Value = 264981.75
$ws.Range("L96").Value = 794945.25
$ws.Range("N96").Value = -799063.25
$ws.Range("H138").Value = 8469423
$ws.Range("J138").Value = 16670833
$ws.Range("L138").Value = 50012499
$ws.Range("N138").Value = -50022779
# row 95 structural update
$ws.Range("H95").Value = 12027
$ws.Range("J95").Value = 12027
$ws.Range("L95").Value = 36081
$ws.Range("N95").Value = -40199

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3779.8235
$ws.Range("I80").Value = 4075.6
$ws.Range("J80").Value = 3357.2856
$ws.Range("K80").Value = 4075.6
$ws.Range("L80").Value = 3357.2856
$ws.Range("M80").Value = -3077.6
$ws.Range("N80").Value = -5353.2856
$ws.Range("H83").Value = 3779.8235
$ws.Range("I83").Value = 4075.6
$ws.Range("J83").Value = 3357.2856
$ws.Range("K83").Value = 20378
$ws.Range("L83").Value = 16786.428
$ws.Range("M83").Value = -15386
$ws.Range("N83").Value = -26770.428
$ws.Range("H97").Value = 782.3929000000001
$ws.Range("I97").Value = 536.9
$ws.Range("J97").Value = 1396.125
$ws.Range("K97").Value = 536.9
$ws.Range("L97").Value = 1396.125
$ws.Range("M97").Value = -40.89999999999998
$ws.Range("N97").Value = -2388.125
$ws.Range("H107").Value = 801.2105
$ws.Range("I107").Value = 666
$ws.Range("J107").Value = 951.44446
$ws.Range("K107").Value = 666
$ws.Range("L107").Value = 951.44446
$ws.Range("M107").Value = 1254
$ws.Range("N107").Value = -4791.44446
$ws.Range("H113").Value = 1729.6364
$ws.Range("J113").Value = 1394.5
$ws.Range("L113").Value = 1394.5
$ws.Range("N113").Value = -5734.5
$ws.Range("H114").Value = 100000
$ws.Range("J114").Value = 100000
$ws.Range("L114").Value = 100000
$ws.Range("N114").Value = -108678
$ws.Range("H131").Value = 78163
$ws.Range("J131").Value = 78163
$ws.Range("L131").Value = 78163
$ws.Range("N131").Value = -88243
$ws.Range("H136").Value = 27781.117
$ws.Range("J136").Value = 27781.117
$ws.Range("L136").Value = 83343.351
$ws.Range("N136").Value = -88443.351

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 31985
$ws.Range("I61").Value = 36107.25
$ws.Range("K61").Value = 36107.25
$ws.Range("M61").Value = -35905.25
$ws.Range("H100").Value = 4166.6665
$ws.Range("I100").Value = 3500
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 3500
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -2959
$ws.Range("N100").Value = -6082
$ws.Range("H113").Value = 31985
$ws.Range("I113").Value = 36107.25
$ws.Range("K113").Value = 36107.25
$ws.Range("M113").Value = -33937.25
$ws.Range("H130").Value = 103622.5
$ws.Range("J130").Value = 103622.5
$ws.Range("L130").Value = 103622.5
$ws.Range("N130").Value = -113662.5
$ws.Range("H132").Value = 3085.6365
$ws.Range("I132").Value = 3071.8838
$ws.Range("J132").Value = 3134.9167
$ws.Range("K132").Value = 9215.651400000001
$ws.Range("L132").Value = 9404.750100000001
$ws.Range("M132").Value = -6685.651400000001
$ws.Range("N132").Value = -14464.7501
$ws.Range("H133").Value = 55648.355
$ws.Range("J133").Value = 55648.355
$ws.Range("L133").Value = 55648.355
$ws.Range("N133").Value = -60708.355
$ws.Range("H136").Value = 4980.5483
$ws.Range("I136").Value = 4649.773
$ws.Range("J136").Value = 5789.1113
$ws.Range("K136").Value = 13949.319
$ws.Range("L136").Value = 17367.3339
$ws.Range("M136").Value = -11399.319
$ws.Range("N136").Value = -22467.3339
# row 137 structural update
$ws.Range("H137").Value = 89333.336
$ws.Range("I137").Value = 80000
$ws.Range("K137").Value = 80000
$ws.Range("M137").Value = -74900
# row 141 structural update
$ws.Range("H141").Value = 93333.336
$ws.Range("J141").Value = 110000
$ws.Range("L141").Value = 110000
$ws.Range("N141").Value = -120360

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1049.8823
$ws.Range("I107").Value = 516.875
$ws.Range("J107").Value = 1523.6666
$ws.Range("K107").Value = 1550.625
$ws.Range("L107").Value = 4570.9998
$ws.Range("M107").Value = 369.375
$ws.Range("N107").Value = -8410.9998
$ws.Range("H132").Value = 3677.1667
$ws.Range("I132").Value = 3922.1
$ws.Range("K132").Value = 11766.3
$ws.Range("M132").Value = -9236.299999999999
$ws.Range("H136").Value = 12695.223
$ws.Range("I136").Value = 12695.223
$ws.Range("K136").Value = 38085.669
$ws.Range("M136").Value = -35535.669
# row 117 structural update
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
